# Alice is the best.
# Apply the Spring 2023 dive-log update to the "Corrected" sheet:
#  1) Rename the T0_* labels used for the ELOW chamber rows to Tn_*
#     (rows 167-170 & 174-177).
#  2) Append the new dive (serial 45091) LOW-chamber block as rows 181-194,
#     mirroring the structure of rows 167-180.
#  3) Move the viewport/selection to the newly-added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")
$ws.Activate()

function Set-StyledCell($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 0) Formatting first (pure style copy, no values => no shared-string churn)
# ---------------------------------------------------------------------
$newRowNums = @(181, 182, 183, 184, 185, 186, 187, 188, 189, 190, 191, 192, 193, 194)
$bStyleSrc = @{ 181 = 167; 182 = 167; 183 = 167; 184 = 167; 188 = 174; 189 = 174; 190 = 174; 191 = 167 }

foreach ($row in $newRowNums) {
    $aDst = "A" + $row
    $hDst = "H" + $row
    Set-StyledCell "A167" $aDst
    Set-StyledCell "H167" $hDst
    if ($bStyleSrc.ContainsKey($row)) {
        $bSrcRow = $bStyleSrc[$row]
        $bSrc = "B" + $bSrcRow
        $bDst = "B" + $row
        Set-StyledCell $bSrc $bDst
    }
}

# ---------------------------------------------------------------------
# 1) T0_t1_ELOW_* -> Tn_t1_ELOW_* relabelling on the existing rows
#    (adds shared strings 192-199, in this row order)
# ---------------------------------------------------------------------
$ws.Range("B167").Value = "Tn_t1_ELOW_tile_01"
$ws.Range("B168").Value = "Tn_t1_ELOW_tile_02"
$ws.Range("B169").Value = "Tn_t1_ELOW_tile_03"
$ws.Range("B170").Value = "Tn_t1_ELOW_blank_01"
$ws.Range("B174").Value = "Tn_t1_ELOW_tile_04"
$ws.Range("B175").Value = "Tn_t1_ELOW_tile_05"
$ws.Range("B176").Value = "Tn_t1_ELOW_tile_06"
$ws.Range("B177").Value = "Tn_t1_ELOW_blank_02"

# ---------------------------------------------------------------------
# 2) New-row B labels (adds shared strings 200-207, in this row order)
# ---------------------------------------------------------------------
$ws.Range("B181").Value = "Tn_t1_LOW_tile_01"
$ws.Range("B182").Value = "Tn_t1_LOW_tile_02"
$ws.Range("B183").Value = "Tn_t1_LOW_tile_03"
$ws.Range("B184").Value = "Tn_t1_LOW_blank_01"
$ws.Range("B188").Value = "Tn_t1_LOW_tile_04"
$ws.Range("B189").Value = "Tn_t1_LOW_tile_05"
$ws.Range("B190").Value = "Tn_t1_LOW_tile_06"
$ws.Range("B191").Value = "Tn_t1_LOW_blank_02"

# ---------------------------------------------------------------------
# 3) New-row H labels (adds shared strings 208-213, first-seen row order)
# ---------------------------------------------------------------------
$ws.Range("H181").Value = "L1"
$ws.Range("H182").Value = "L2"
$ws.Range("H183").Value = "L3"
$ws.Range("H185").Value = "L4"
$ws.Range("H186").Value = "L5"
$ws.Range("H187").Value = "L6"
$ws.Range("H188").Value = "L4"
$ws.Range("H189").Value = "L5"
$ws.Range("H190").Value = "L6"
$ws.Range("H192").Value = "L1"
$ws.Range("H193").Value = "L2"
$ws.Range("H194").Value = "L3"

# ---------------------------------------------------------------------
# 4) Remaining columns for the new rows: A (date), C (chamber/phase), G (tile#)
#    These reuse existing shared strings (Light_0x / Dark_0x) or are plain
#    numbers, so they don't disturb the shared-string table ordering.
# ---------------------------------------------------------------------
foreach ($row in $newRowNums) {
    $ws.Range("A" + $row).Value = 45091
}

$ws.Range("C181").Value = "Light_01"
$ws.Range("C182").Value = "Light_02"
$ws.Range("C183").Value = "Light_03"
$ws.Range("C184").Value = "Light_04"
$ws.Range("C185").Value = "Dark_01"
$ws.Range("C186").Value = "Dark_02"
$ws.Range("C187").Value = "Dark_03"
$ws.Range("C188").Value = "Light_05"
$ws.Range("C189").Value = "Light_06"
$ws.Range("C190").Value = "Light_07"
$ws.Range("C191").Value = "Light_08"
$ws.Range("C192").Value = "Dark_05"
$ws.Range("C193").Value = "Dark_06"
$ws.Range("C194").Value = "Dark_07"

$ws.Range("G181").Value = 1
$ws.Range("G182").Value = 2
$ws.Range("G183").Value = 3
$ws.Range("G184").Value = 4
$ws.Range("G185").Value = 5
$ws.Range("G186").Value = 6
$ws.Range("G187").Value = 7
$ws.Range("G188").Value = 5
$ws.Range("G189").Value = 6
$ws.Range("G190").Value = 7
$ws.Range("G191").Value = 4
$ws.Range("G192").Value = 1
$ws.Range("G193").Value = 2
$ws.Range("G194").Value = 3

# ---------------------------------------------------------------------
# 5) Viewport / selection bookkeeping
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 164
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J189").Select()
